$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.091.71"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "3.504.92"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.29"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.35"
$ws.Range("E6").Value = "  +2.31%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.601"
$ws.Range("E8").Value = "  +2.57%  "
$ws.Range("E9").Value = "  +4.58%  "
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.434"
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").Value = "4.108.40"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.03"
$ws.Range("E14").Value = "  +2.67%  "
$ws.Range("D15").Value = "67.077.05"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000179"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").Value = "3.471.56"
$ws.Range("E17").Value = "  -1.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.32"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "394.28"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.02"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.14"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.537"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.68"
$ws.Range("E25").Value = "  -3.39%  "
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.26"
$ws.Range("E27").Value = "  -1.93%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.28"
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.77"
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.36"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("E35").Value = "  +3.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.42"
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.881"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.99"
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.69"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0747"
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.32"
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.32"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").Value = "2.814.52"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0303"
$ws.Range("E47").Value = "  -2.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "335.61"
$ws.Range("E48").Value = "  -5.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.63"
$ws.Range("E49").Value = "  +2.75%  "
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("E51").Value = "  +0.64%  "
# Rows 45 and 46 swap positions (OKB <-> dogwifhat) with updated values
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.57"
$ws.Range("E45").Value = "  +2.22%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.70"
$ws.Range("E46").Value = "  -0.98%  "
